$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("I2").Value = 7004
$ws.Range("I3").Value = 7247
$ws.Range("B4").Value = 1668
$ws.Range("C4").Value = 1816
$ws.Range("E4").Value = 1970
$ws.Range("F4").Value = 1873
$ws.Range("G4").Value = 1447
$ws.Range("I4").Value = 1669
$ws.Range("I6").Value = 8553
$ws.Range("B7").Value = 23300
$ws.Range("C7").Value = 28359
$ws.Range("E7").Value = 25974
$ws.Range("F7").Value = 24062
$ws.Range("G7").Value = 24672
$ws.Range("I7").Value = 25154

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("I6").Value = 121
$ws.Range("I7").Value = 295

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("I4").Value = 7
$ws.Range("I7").Value = 141

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range("I6").Value = 23
$ws.Range("I7").Value = 83

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("I3").Value = 252
$ws.Range("I7").Value = 776

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("I6").Value = 114
$ws.Range("I7").Value = 438

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("I3").Value = 356
$ws.Range("I6").Value = 284
$ws.Range("I7").Value = 952

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("I2").Value = 82
$ws.Range("I7").Value = 222

$ws = $wb.Worksheets.Item('New City')
$ws.Range("I2").Value = 187
$ws.Range("I4").Value = 25
$ws.Range("I7").Value = 586

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("I2").Value = 200
$ws.Range("I7").Value = 788
$ws.Range("I8").Value = 1494
$ws.Range("I11").Value = 380
$ws.Range("I14").Value = 141
$ws.Range("I18").Value = 198
$ws.Range("I20").Value = 623
$ws.Range("I27").Value = 219
$ws.Range("I29").Value = 1497
$ws.Range("I30").Value = 83
$ws.Range("I33").Value = 1109
$ws.Range("I34").Value = 114
$ws.Range("I37").Value = 776
$ws.Range("I40").Value = 44
$ws.Range("I42").Value = 943
$ws.Range("I43").Value = 218
$ws.Range("I44").Value = 192
$ws.Range("I46").Value = 58
$ws.Range("I49").Value = 168
$ws.Range("I50").Value = 133
$ws.Range("I52").Value = 570
$ws.Range("I53").Value = 284
$ws.Range("I54").Value = 492
$ws.Range("B63").Value = 372
$ws.Range("C63").Value = 248
$ws.Range("E63").Value = 316
$ws.Range("F63").Value = 163
$ws.Range("G63").Value = 208
$ws.Range("I63").Value = 84
$ws.Range("I65").Value = 586
$ws.Range("I67").Value = 952
$ws.Range("I70").Value = 44
$ws.Range("I73").Value = 228
$ws.Range("I76").Value = 357
$ws.Range("I79").Value = 721
$ws.Range("I83").Value = 542
$ws.Range("I84").Value = 222
$ws.Range("I85").Value = 1120
$ws.Range("I87").Value = 66
$ws.Range("I90").Value = 327
$ws.Range("I95").Value = 387
$ws.Range("I96").Value = 295
$ws.Range("I98").Value = 184
$ws.Range("I99").Value = 438
$ws.Range("B101").Value = 23300
$ws.Range("C101").Value = 28359
$ws.Range("E101").Value = 25974
$ws.Range("F101").Value = 24062
$ws.Range("G101").Value = 24672
$ws.Range("I101").Value = 25154

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("I3").Value = 197
$ws.Range("I7").Value = 542

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("I6").Value = 79
$ws.Range("I7").Value = 387

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("I2").Value = 250
$ws.Range("I3").Value = 408
$ws.Range("I6").Value = 357
$ws.Range("I7").Value = 1109

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("I2").Value = 32
$ws.Range("I7").Value = 168

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("I4").Value = 33
$ws.Range("I6").Value = 239
$ws.Range("I7").Value = 492

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("I2").Value = 445
$ws.Range("I4").Value = 81
$ws.Range("I6").Value = 412
$ws.Range("I7").Value = 1497

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("I4").Value = 17
$ws.Range("I7").Value = 192

$ws = $wb.Worksheets.Item('River North')
$ws.Range("I6").Value = 163
$ws.Range("I7").Value = 357

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("I3").Value = 422
$ws.Range("I6").Value = 294
$ws.Range("I7").Value = 1120

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("I3").Value = 272
$ws.Range("I7").Value = 943

$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Range("I2").Value = 17
$ws.Range("I3").Value = 20
$ws.Range("I7").Value = 58

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("I3").Value = 238
$ws.Range("I7").Value = 721

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("I6").Value = 218
$ws.Range("I7").Value = 623

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("I6").Value = 92
$ws.Range("I7").Value = 198

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("I2").Value = 142
$ws.Range("I7").Value = 570

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range("I6").Value = 24
$ws.Range("I7").Value = 114

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("I6").Value = 120
$ws.Range("I7").Value = 184

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range("I6").Value = 43
$ws.Range("I7").Value = 133

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("I2").Value = 150
$ws.Range("I7").Value = 380

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("I4").Value = 24
$ws.Range("I7").Value = 228

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("I6").Value = 43
$ws.Range("I7").Value = 200

$ws = $wb.Worksheets.Item('O''Hare')
$ws.Range("I3").Value = 13
$ws.Range("I7").Value = 44

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("I4").Value = 94
$ws.Range("I7").Value = 1494

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("I4").Value = 30
$ws.Range("I7").Value = 219

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("I3").Value = 85
$ws.Range("I7").Value = 327

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("I2").Value = 47
$ws.Range("I6").Value = 120
$ws.Range("I7").Value = 218

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("I2").Value = 60
$ws.Range("I3").Value = 57
$ws.Range("I7").Value = 284

$ws = $wb.Worksheets.Item('Hegewisch')
$ws.Range("I4").Value = 2
$ws.Range("I7").Value = 44

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("I6").Value = 214
$ws.Range("I7").Value = 788

$ws = $wb.Worksheets.Item('Ukrainian Village')
$ws.Range("I2").Value = 8
$ws.Range("I7").Value = 66
